$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Remove the "stay informed about IBM training" block of paragraphs
# (the introductory sentence plus the four social/news bullet lines:
# IBM Training News, YouTube, Facebook, Twitter) from the "For more
# information" section, collapsing the introductory paragraph into a
# single empty paragraph that keeps the indentation used by the
# bulleted lines that followed it.
# ------------------------------------------------------------------

function Get-ParagraphIndexForRange($doc, $rngStart, $rngEnd) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -eq $rngStart -and $p.Range.End -eq $rngEnd) {
            return $i
        }
    }
    return -1
}

# Locate the full paragraph that introduces the block.
$introRange = $d.Content
$foundIntro = $introRange.Find.Execute("To stay informed about IBM training, see the following sites:")
if (-not $foundIntro) {
    throw "Could not find the 'To stay informed about IBM training' paragraph."
}
$introRange.Expand(4) | Out-Null   # wdParagraph -> expand the hit to the whole paragraph
$introIndex = Get-ParagraphIndexForRange $d $introRange.Start $introRange.End
if ($introIndex -lt 0) {
    throw "Could not resolve the paragraph index of the intro paragraph."
}

# Locate the full paragraph of the final line of the block (Twitter).
$lastRange = $d.Content
$foundLast = $lastRange.Find.Execute("Twitter: ")
if (-not $foundLast) {
    throw "Could not find the 'Twitter:' paragraph."
}
$lastRange.Expand(4) | Out-Null    # wdParagraph -> expand the hit to the whole paragraph

# Delete everything from the end of the intro paragraph's text (i.e.
# right after its own paragraph mark) through the end of the Twitter
# paragraph. This removes the IBM Training News / YouTube / Facebook
# / Twitter paragraphs in their entirety, including their paragraph
# marks, while leaving the intro paragraph's own mark untouched.
$afterIntroStart = $introRange.End
$blockEnd = $lastRange.End
if ($blockEnd -gt $afterIntroStart) {
    $killRange = $d.Range($afterIntroStart, $blockEnd)
    $killRange.Delete()
}

# Clear the text of the intro paragraph itself, leaving only its
# paragraph mark.
$introTextRange = $d.Range($introRange.Start, $introRange.End - 1)
if ($introTextRange.End -gt $introTextRange.Start) {
    $introTextRange.Delete()
}

# Re-fetch the (now empty) paragraph by its stable collection index
# and give it the left indent the removed bullet paragraphs used to
# have (0.5in / 36pt / 720 twips). Using the Paragraphs collection
# directly avoids ambiguity that comes from formatting a zero-length
# Range sitting exactly on a paragraph boundary.
$introPara = $d.Paragraphs.Item($introIndex)
$introPara.LeftIndent = 36
